# CPL.xlsx update for the Pixel8x4U add-on:
#  - Rename the worksheet to reflect the board (GD32F207RG -> GD32F-Pixel8x4U)
#  - Refresh the component placement table (Designator / Mid X / Mid Y / Layer / Rotation)
#    with the full BOM for the new revision, including the previously-empty rows
#    (C5, C6, R12-R20, SW1, U2) that add support for clock based LED chips (e.g. APA102)
#  - Tweak column widths and the saved cursor/selection state to match the authoring tool

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet ---
$ws.Name = "GD32F-Pixel8x4U"

# --- Component placement data (row -> Designator, Mid X, Mid Y, Rotation); Layer is always "Top" ---
$data = @(
    @("C1", 44.45, 15.24, 90),
    @("C2", 42.23, 15.24, 90),
    @("C3", 33.02, 14.29, 180),
    @("C4", 29.21, 29.85, 0),
    @("C5", 33.02, 29.85, 180),
    @("C6", 29.21, 14.29, 0),
    @("LED1", 46.04, 12.7, 180),
    @("Q1", 47.31, 15.24, 180),
    @("Q2", 47.31, 19.37, 90),
    @("R1", 48.9, 22.54, 180),
    @("R2", 44.77, 18.73, 270),
    @("R3", 42.86, 12.7, 0),
    @("R4", 46.99, 31.75, 270),
    @("R5", 23.18, 13.97, 180),
    @("R6", 23.18, 15.24, 180),
    @("R7", 23.18, 16.51, 180),
    @("R8", 23.18, 17.78, 180),
    @("R9", 23.18, 19.05, 180),
    @("R10", 23.18, 20.32, 180),
    @("R11", 23.18, 21.59, 180),
    @("R12", 23.18, 22.86, 180),
    @("R13", 23.18, 31.43, 180),
    @("R14", 23.18, 32.7, 180),
    @("R15", 23.18, 33.97, 180),
    @("R16", 23.18, 35.24, 180),
    @("R17", 23.18, 36.51, 180),
    @("R18", 23.18, 37.78, 180),
    @("R19", 23.18, 39.05, 180),
    @("R20", 23.18, 40.32, 180),
    @("SW1", 48.9, 38.73, 90),
    @("U1", 31.12, 19.05, 90),
    @("U2", 31.12, 34.93, 90)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "Top"
    $ws.Cells.Item($r, 5).Value = $row[3]
    $r++
}

# --- Column widths: narrow col A slightly, and give col E its own explicit width ---
$ws.Columns("A").ColumnWidth = 8.75
$ws.Columns("E").ColumnWidth = 7.59

# --- Restore the saved selection/cursor state (active cell A2) ---
$ws.Range("A2").Select()
